$wb = $excel.ActiveWorkbook

# --- Add the new "Old" worksheet (Worksheets.Add() inserts before the
# currently active sheet, i.e. before "CDRStatus", giving us the desired
# sheet order: Old, CDRStatus) ---
$old = $wb.Worksheets.Add()
$old.Name = "Old"

# Populate the "Old" sheet. Put A3 (which reuses the original
# "CDR.LDA0610547" shared string) in before we repoint CDRStatus!A2 away
# from it, so that string is preserved (not garbage-collected) and the new
# string "CDR.CLK0601316" is appended as a new shared-string entry, matching
# the target shared string ordering.
$old.Range("A1").Value = "value:1:1:1"
$old.Range("A3").Value = "CDR.LDA0610547"
$old.Range("A2").Value = "CDR.CLK0601316"

# Re-point the cell selection/active cell on the "Old" sheet to A3.
$old.Range("A3").Select() | Out-Null

# Fit column A to its contents.
$old.Columns("A").AutoFit() | Out-Null

# Give the (empty) B8 cell the new "Trebuchet MS" / 10pt / black font used
# by the workbook, and size row 8 accordingly.
$oldFont = $old.Range("B8").Font
$oldFont.Name = "Trebuchet MS"
$oldFont.Size = 10
$oldFont.Color = 0
$old.Rows(8).RowHeight = 15.75

# Portrait page setup for the new sheet.
$old.PageSetup.Orientation = 1

# --- Update the "CDRStatus" sheet: A2 now references the new string ---
$cdr = $wb.Worksheets.Item("CDRStatus")
$cdr.Range("A2").Value = "CDR.CLK0601316"

# --- Make "CDRStatus" (the second tab) the active sheet/tab ---
$old.Activate() | Out-Null
$cdr.Activate() | Out-Null
